# Update "paises" (countries) COVID tracking sheet + provincias Spain data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header timestamp update ---
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 08:52"

# --- Singapur overtakes Japon/Chile/Pakistan (rows 29-32) ---
# Row 29 becomes Singapur with fresh data; Japon/Chile/Pakistan each shift down one row.
$ws.Cells.Item(29, 1).Value = "Singapur"
$ws.Cells.Item(29, 2).Value = 13624
$ws.Cells.Item(29, 3).Value = 931
$ws.Cells.Item(29, 4).Value = 1002
$ws.Cells.Item(29, 5).Value = 12610
$ws.Cells.Item(29, 6).Value = 24
$ws.Cells.Item(29, 7).Value = 0
$ws.Cells.Item(29, 8).Value = 12

$ws.Cells.Item(30, 1).Value = "Japon"
$ws.Cells.Item(30, 2).Value = 13231
$ws.Cells.Item(30, 3).Value = 0
$ws.Cells.Item(30, 4).Value = 1656
$ws.Cells.Item(30, 5).Value = 11215
$ws.Cells.Item(30, 6).Value = 287
$ws.Cells.Item(30, 7).Value = 0
$ws.Cells.Item(30, 8).Value = 360

$ws.Cells.Item(31, 1).Value = "Chile"
$ws.Cells.Item(31, 2).Value = 12858
$ws.Cells.Item(31, 3).Value = 0
$ws.Cells.Item(31, 4).Value = 6746
$ws.Cells.Item(31, 5).Value = 5931
$ws.Cells.Item(31, 6).Value = 418
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 181

$ws.Cells.Item(32, 1).Value = "Pakistan"
$ws.Cells.Item(32, 2).Value = 12723
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(32, 4).Value = 2866
$ws.Cells.Item(32, 5).Value = 9588
$ws.Cells.Item(32, 6).Value = 111
$ws.Cells.Item(32, 7).Value = 0
$ws.Cells.Item(32, 8).Value = 269

# --- Ucrania overtakes Indonesia/Dinamarca (rows 39-41) ---
# Row 39 becomes Ucrania with fresh data; Indonesia/Dinamarca each shift down one row.
$ws.Cells.Item(39, 1).Value = "Ucrania"
$ws.Cells.Item(39, 2).Value = 8617
$ws.Cells.Item(39, 3).Value = 492
$ws.Cells.Item(39, 4).Value = 840
$ws.Cells.Item(39, 5).Value = 7568
$ws.Cells.Item(39, 6).Value = 107
$ws.Cells.Item(39, 7).Value = 8
$ws.Cells.Item(39, 8).Value = 209

$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(40, 2).Value = 8607
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 1042
$ws.Cells.Item(40, 5).Value = 6845
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 720

$ws.Cells.Item(41, 1).Value = "Dinamarca"
$ws.Cells.Item(41, 2).Value = 8445
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 5669
$ws.Cells.Item(41, 5).Value = 2358
$ws.Cells.Item(41, 6).Value = 70
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = 418

# --- Tailandia (row 60) gets refreshed case counts ---
$ws.Cells.Item(60, 2).Value = 2922
$ws.Cells.Item(60, 3).Value = 15
$ws.Cells.Item(60, 4).Value = 2594
$ws.Cells.Item(60, 5).Value = 277
